$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells H1:J1
$ws.Range("H1").Value = "grpf_simul"
$ws.Range("I1").Value = "vu_simul"
$ws.Range("J1").Value = "shortage_simul"

# Copy style from the existing header cell (G1) to the new header cells
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for columns H, I, J (rows 2-19)
$data = @{
    2  = @(-0.5064595937728882, 1.192999362945557, 4.666666507720947)
    3  = @(0.2062890231609344, 1.165635347366333, 5)
    4  = @(-2.330345392227173, 1.221376657485962, 5.666666507720947)
    5  = @(-2.08102560043335, 1.20171332359314, 4.666666507720947)
    6  = @(-0.2500871419906616, 1.184739232063293, 4.666666507720947)
    7  = @(-0.6588118672370911, 1.073658347129822, 12)
    8  = @(8.748804092407227, 0.2628340423107147, 14)
    9  = @(-3.022449254989624, 0.4582114815711975, 25.33333396911621)
    10 = @(-0.6266276240348816, 0.6320492029190063, 9.666666984558105)
    11 = @(-3.191989183425903, 0.7893379926681519, 11.66666698455811)
    12 = @(2.508954286575317, 1.009974718093872, 48.33333206176758)
    13 = @(1.477425456047058, 1.300251245498657, 22.33333396911621)
    14 = @(4.266705513000488, 1.632708787918091, 22)
    15 = @(5.85883092880249, 1.84051239490509, 19.66666603088379)
    16 = @(5.80972146987915, 1.911138534545898, 32.66666793823242)
    17 = @(6.041294574737549, 1.843034982681274, 18)
    18 = @(2.111628532409668, 1.825654029846191, 23.66666603088379)
    19 = @(-1.174693703651428, 1.724597811698914, 21.5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
    $ws.Cells.Item($row, 10).Value = $vals[2]
}
